$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.476.63"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.363.48"
$ws.Range("E3").Value = "  -2.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.24"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.31"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.358.75"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.105"
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.06"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.336"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.45"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.781.56"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  -3.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.389.34"
$ws.Range("E17").Value = "  -2.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.349.71"
$ws.Range("E18").Value = "  -3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.03"
$ws.Range("E19").Value = "  +10.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.43"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.23"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.02"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.95"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.80"
$ws.Range("E25").Value = "  -6.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.08"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "550.54"
$ws.Range("E27").Value = "  -3.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  -9.84%  "
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0912"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.30"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.77"
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.41"
$ws.Range("E36").Value = "  +2.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "150.72"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.365"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.52"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.06"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.99"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.37"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.64"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0286"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "138.35"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.49"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.582"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0498"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.07"
$ws.Range("E51").Value = "  -2.21%  "
